# Auto-generated: updates market/profit data cells per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3051.62
$ws.Range("J17").Value = 3542.4048
$ws.Range("L17").Value = 10627.2144
$ws.Range("N17").Value = -10963.2144
$ws.Range("H38").Value = 538.5
$ws.Range("I38").Value = 209.44444
$ws.Range("J38").Value = 3500
$ws.Range("K38").Value = 628.33332
$ws.Range("L38").Value = 10500
$ws.Range("M38").Value = -256.33332
$ws.Range("N38").Value = -11244
$ws.Range("H40").Value = 9270.632
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 9270.632
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 9270.632
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -9620.632
$ws.Range("H58").Value = 2005.6666
$ws.Range("I58").Value = 2000
$ws.Range("J58").Value = 2008.5
$ws.Range("K58").Value = 6000
$ws.Range("L58").Value = 6025.5
$ws.Range("M58").Value = -5850
$ws.Range("N58").Value = -6325.5
$ws.Range("H62").Value = 8111.706
$ws.Range("I62").Value = 6315.6665
$ws.Range("J62").Value = 9091.362999999999
$ws.Range("K62").Value = 6315.6665
$ws.Range("L62").Value = 9091.362999999999
$ws.Range("M62").Value = -5691.6665
$ws.Range("N62").Value = -10339.363
$ws.Range("H65").Value = 8111.706
$ws.Range("I65").Value = 6315.6665
$ws.Range("J65").Value = 9091.362999999999
$ws.Range("K65").Value = 31578.3325
$ws.Range("L65").Value = 45456.815
$ws.Range("M65").Value = -28458.3325
$ws.Range("N65").Value = -51696.815
$ws.Range("H70").Value = 8003857
$ws.Range("I70").Value = 3460.9167
$ws.Range("J70").Value = 15388838
$ws.Range("K70").Value = 10382.7501
$ws.Range("L70").Value = 46166514
$ws.Range("M70").Value = -10112.7501
$ws.Range("N70").Value = -46167054
$ws.Range("H73").Value = 8003857
$ws.Range("I73").Value = 3460.9167
$ws.Range("J73").Value = 15388838
$ws.Range("K73").Value = 10382.7501
$ws.Range("L73").Value = 46166514
$ws.Range("M73").Value = -9446.750100000001
$ws.Range("N73").Value = -46168386
$ws.Range("H100").Value = 3305.0833
$ws.Range("I100").Value = 2070.1428
$ws.Range("J100").Value = 5034
$ws.Range("K100").Value = 2070.1428
$ws.Range("L100").Value = 5034
$ws.Range("M100").Value = -1529.1428
$ws.Range("N100").Value = -6116
$ws.Range("H112").Value = 1535.1
$ws.Range("I112").Value = 1271.1666
$ws.Range("J112").Value = 1648.2142
$ws.Range("K112").Value = 3813.4998
$ws.Range("L112").Value = 4944.642599999999
$ws.Range("M112").Value = -2705.4998
$ws.Range("N112").Value = -7160.642599999999
$ws.Range("H119").Value = 500
$ws.Range("J119").Value = 500
$ws.Range("L119").Value = 1500
$ws.Range("N119").Value = -11176
$ws.Range("H121").Value = 1536.2
$ws.Range("J121").Value = 1536.2
$ws.Range("L121").Value = 4608.6
$ws.Range("N121").Value = -8102.6
$ws.Range("H131").Value = 5298.1816
$ws.Range("I131").Value = 4780.4165
$ws.Range("K131").Value = 14341.2495
$ws.Range("M131").Value = -9301.249500000002
$ws.Range("H137").Value = 76929130
$ws.Range("I137").Value = 200003710
$ws.Range("J137").Value = 7513
$ws.Range("K137").Value = 600011130
$ws.Range("L137").Value = 22539
$ws.Range("M137").Value = -600008580
$ws.Range("N137").Value = -27639
$ws.Range("H138").Value = 5380.1143
$ws.Range("I138").Value = 3772.8333
$ws.Range("J138").Value = 6218.696
$ws.Range("K138").Value = 11318.4999
$ws.Range("L138").Value = 18656.088
$ws.Range("M138").Value = -6178.499899999999
$ws.Range("N138").Value = -28936.088
$ws.Range("H141").Value = 3293.524
$ws.Range("I141").Value = 2421.4119
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 7264.2357
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = -2084.2357
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1684.1945
$ws.Range("I32").Value = 1143.6119
$ws.Range("J32").Value = 8928
$ws.Range("K32").Value = 1143.6119
$ws.Range("L32").Value = 8928
$ws.Range("M32").Value = -856.6119000000001
$ws.Range("N32").Value = -9502
$ws.Range("H45").Value = 100005180
$ws.Range("I45").Value = 200001840
$ws.Range("K45").Value = 200001840
$ws.Range("M45").Value = -200001463
$ws.Range("H61").Value = 7362.28
$ws.Range("I61").Value = 6020.45
$ws.Range("J61").Value = 12729.6
$ws.Range("K61").Value = 6020.45
$ws.Range("L61").Value = 12729.6
$ws.Range("M61").Value = -5808.45
$ws.Range("N61").Value = -13153.6
$ws.Range("H97").Value = 1592.5758
$ws.Range("J97").Value = 4244
$ws.Range("L97").Value = 4244
$ws.Range("N97").Value = -5236
$ws.Range("H132").Value = 3868.7778
$ws.Range("I132").Value = 2698.8572
$ws.Range("K132").Value = 8096.571599999999
$ws.Range("M132").Value = -5566.571599999999
$ws.Range("H136").Value = 7362.28
$ws.Range("I136").Value = 6020.45
$ws.Range("J136").Value = 12729.6
$ws.Range("K136").Value = 18061.35
$ws.Range("L136").Value = 38188.8
$ws.Range("M136").Value = -15511.35
$ws.Range("N136").Value = -43288.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2522.879
$ws.Range("I86").Value = 1958.75
$ws.Range("J86").Value = 3390.7693
$ws.Range("K86").Value = 1958.75
$ws.Range("L86").Value = 3390.7693
$ws.Range("M86").Value = -835.75
$ws.Range("N86").Value = -5636.7693
$ws.Range("H89").Value = 2522.879
$ws.Range("I89").Value = 1958.75
$ws.Range("J89").Value = 3390.7693
$ws.Range("K89").Value = 9793.75
$ws.Range("L89").Value = 16953.8465
$ws.Range("M89").Value = -4177.75
$ws.Range("N89").Value = -28185.8465
$ws.Range("H94").Value = 1753.4062
$ws.Range("I94").Value = 1350.5
$ws.Range("J94").Value = 2962.125
$ws.Range("K94").Value = 1350.5
$ws.Range("L94").Value = 2962.125
$ws.Range("M94").Value = -899.5
$ws.Range("N94").Value = -3864.125
$ws.Range("H107").Value = 1978.725
$ws.Range("I107").Value = 1881.3715
$ws.Range("J107").Value = 2660.2
$ws.Range("K107").Value = 1881.3715
$ws.Range("L107").Value = 2660.2
$ws.Range("M107").Value = 38.62850000000003
$ws.Range("N107").Value = -6500.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35070.113
$ws.Range("I31").Value = 4293.5
$ws.Range("J31").Value = 76105.60000000001
$ws.Range("K31").Value = 4293.5
$ws.Range("L31").Value = 76105.60000000001
$ws.Range("M31").Value = -3998.5
$ws.Range("N31").Value = -76695.60000000001
$ws.Range("H34").Value = 35070.113
$ws.Range("I34").Value = 4293.5
$ws.Range("J34").Value = 76105.60000000001
$ws.Range("K34").Value = 4293.5
$ws.Range("L34").Value = 76105.60000000001
$ws.Range("M34").Value = -4091.5
$ws.Range("N34").Value = -76509.60000000001
$ws.Range("H107").Value = 2410.7144
$ws.Range("I107").Value = 2975.6
$ws.Range("J107").Value = 998.5
$ws.Range("K107").Value = 2975.6
$ws.Range("L107").Value = 998.5
$ws.Range("M107").Value = -1055.6
$ws.Range("N107").Value = -4838.5
$ws.Range("H134").Value = 2751.8865
$ws.Range("I134").Value = 2091.45
$ws.Range("J134").Value = 9356.25
$ws.Range("K134").Value = 6274.349999999999
$ws.Range("L134").Value = 28068.75
$ws.Range("M134").Value = -3739.349999999999
$ws.Range("N134").Value = -33138.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 913674.1
$ws.Range("I60").Value = 1093.5
$ws.Range("K60").Value = 3280.5
$ws.Range("M60").Value = -3029.5
$ws.Range("H68").Value = 5685967
$ws.Range("I68").Value = 10418037
$ws.Range("J68").Value = 7483.1
$ws.Range("K68").Value = 31254111
$ws.Range("L68").Value = 22449.3
$ws.Range("M68").Value = -31253300
$ws.Range("N68").Value = -24071.3
$ws.Range("H71").Value = 5685967
$ws.Range("I71").Value = 10418037
$ws.Range("J71").Value = 7483.1
$ws.Range("K71").Value = 93762333
$ws.Range("L71").Value = 67347.90000000001
$ws.Range("M71").Value = -93758277
$ws.Range("N71").Value = -75459.90000000001
$ws.Range("H93").Value = 22513
$ws.Range("J93").Value = 22513
$ws.Range("L93").Value = 67539
$ws.Range("N93").Value = -71283
$ws.Range("H99").Value = 10006.25
$ws.Range("I99").Value = 11012.5
$ws.Range("K99").Value = 33037.5
$ws.Range("M99").Value = -30791.5
$ws.Range("H105").Value = 19014.5
$ws.Range("I105").Value = 19000
$ws.Range("J105").Value = 19029
$ws.Range("K105").Value = 57000
$ws.Range("L105").Value = 57087
$ws.Range("M105").Value = -54379
$ws.Range("N105").Value = -62329
$ws.Range("H113").Value = 40001004
$ws.Range("I113").Value = 700.6923
$ws.Range("J113").Value = 83334664
$ws.Range("K113").Value = 2102.0769
$ws.Range("L113").Value = 250003992
$ws.Range("M113").Value = 67.92309999999998
$ws.Range("N113").Value = -250008332
$ws.Range("H129").Value = 4169749.8
$ws.Range("I129").Value = 739.2222
$ws.Range("J129").Value = 7580758.5
$ws.Range("K129").Value = 2217.6666
$ws.Range("L129").Value = 22742275.5
$ws.Range("M129").Value = 2782.3334
$ws.Range("N129").Value = -22752275.5
$ws.Range("H131").Value = 10804957
$ws.Range("I131").Value = 2398.75
$ws.Range("J131").Value = 32410074
$ws.Range("K131").Value = 7196.25
$ws.Range("L131").Value = 97230222
$ws.Range("M131").Value = -2156.25
$ws.Range("N131").Value = -97240302
$ws.Range("H132").Value = 3954.394
$ws.Range("I132").Value = 4037.3333
$ws.Range("J132").Value = 3907
$ws.Range("K132").Value = 36335.9997
$ws.Range("L132").Value = 35163
$ws.Range("M132").Value = -33805.9997
$ws.Range("N132").Value = -40223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7675.25
$ws.Range("I22").Value = 2899.5
$ws.Range("J22").Value = 9267.166999999999
$ws.Range("K22").Value = 2899.5
$ws.Range("L22").Value = 9267.166999999999
$ws.Range("M22").Value = -2604.5
$ws.Range("N22").Value = -9857.166999999999
$ws.Range("H27").Value = 7675.25
$ws.Range("I27").Value = 2899.5
$ws.Range("J27").Value = 9267.166999999999
$ws.Range("K27").Value = 2899.5
$ws.Range("L27").Value = 9267.166999999999
$ws.Range("M27").Value = -2792.5
$ws.Range("N27").Value = -9481.166999999999
$ws.Range("H132").Value = 3502.6287
$ws.Range("I132").Value = 1753.2858
$ws.Range("J132").Value = 10500
$ws.Range("K132").Value = 5259.857400000001
$ws.Range("L132").Value = 31500
$ws.Range("M132").Value = -2729.857400000001
$ws.Range("N132").Value = -36560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1397.1305
$ws.Range("I100").Value = 1141.7
$ws.Range("J100").Value = 3100
$ws.Range("K100").Value = 2283.4
$ws.Range("L100").Value = 6200
$ws.Range("M100").Value = -1742.4
$ws.Range("N100").Value = -7282
$ws.Range("H107").Value = 339.6389
$ws.Range("I107").Value = 297.71875
$ws.Range("J107").Value = 675
$ws.Range("K107").Value = 893.15625
$ws.Range("L107").Value = 2025
$ws.Range("M107").Value = 1026.84375
$ws.Range("N107").Value = -5865
$ws.Range("H122").Value = 4134.1763
$ws.Range("I122").Value = 3440.5833
$ws.Range("J122").Value = 5798.8
$ws.Range("K122").Value = 10321.7499
$ws.Range("L122").Value = 17396.4
$ws.Range("M122").Value = -7871.749899999999
$ws.Range("N122").Value = -22296.4
